$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1168
$ws.Range("I6").Value = 390.66666
$ws.Range("K6").Value = 1171.99998
$ws.Range("M6").Value = -1059.99998

$ws.Range("H38").Value = 13500
$ws.Range("I38").Value = 500
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 1500
$ws.Range("L38").Value = 60000
$ws.Range("M38").Value = -1128
$ws.Range("N38").Value = -60744

$ws.Range("H41").Value = 815.75
$ws.Range("I41").Value = 358
$ws.Range("J41").Value = 1142.7142
$ws.Range("K41").Value = 358
$ws.Range("L41").Value = 1142.7142
$ws.Range("M41").Value = 82
$ws.Range("N41").Value = -2022.7142

$ws.Range("H103").Value = 9836.444
$ws.Range("I103").Value = 438
$ws.Range("J103").Value = 28633.334
$ws.Range("K103").Value = 1314
$ws.Range("L103").Value = 85900.00199999999
$ws.Range("M103").Value = -728
$ws.Range("N103").Value = -87072.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8264.860000000001
$ws.Range("I32").Value = 3825.6365
$ws.Range("J32").Value = 16882.176
$ws.Range("K32").Value = 3825.6365
$ws.Range("L32").Value = 16882.176
$ws.Range("M32").Value = -3538.6365
$ws.Range("N32").Value = -17456.176

$ws.Range("H61").Value = 1911.5
$ws.Range("I61").Value = 1531.8334
$ws.Range("J61").Value = 3620
$ws.Range("K61").Value = 1531.8334
$ws.Range("L61").Value = 3620
$ws.Range("M61").Value = -1319.8334
$ws.Range("N61").Value = -4044

$ws.Range("H134").Value = 48164.535
$ws.Range("J134").Value = 48164.535
$ws.Range("L134").Value = 48164.535
$ws.Range("N134").Value = -58304.535

$ws.Range("H136").Value = 1911.5
$ws.Range("I136").Value = 1531.8334
$ws.Range("J136").Value = 3620
$ws.Range("K136").Value = 4595.5002
$ws.Range("L136").Value = 10860
$ws.Range("M136").Value = -2045.5002
$ws.Range("N136").Value = -15960

$ws.Range("H137").Value = 55750
$ws.Range("J137").Value = 55750
$ws.Range("L137").Value = 55750
$ws.Range("N137").Value = -65950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 40078.75
$ws.Range("J137").Value = 40078.75
$ws.Range("L137").Value = 40078.75
$ws.Range("N137").Value = -50278.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3390.5476
$ws.Range("I31").Value = 1291.2084
$ws.Range("J31").Value = 6189.6665
$ws.Range("K31").Value = 1291.2084
$ws.Range("L31").Value = 6189.6665
$ws.Range("M31").Value = -996.2084
$ws.Range("N31").Value = -6779.6665

$ws.Range("H34").Value = 3390.5476
$ws.Range("I34").Value = 1291.2084
$ws.Range("J34").Value = 6189.6665
$ws.Range("K34").Value = 1291.2084
$ws.Range("L34").Value = 6189.6665
$ws.Range("M34").Value = -1089.2084
$ws.Range("N34").Value = -6593.6665

$ws.Range("H58").Value = 1886.2354
$ws.Range("I58").Value = 1632.25
$ws.Range("J58").Value = 5950
$ws.Range("K58").Value = 1632.25
$ws.Range("L58").Value = 5950
$ws.Range("M58").Value = -1429.25
$ws.Range("N58").Value = -6356

$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497

$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485

$ws.Range("H132").Value = 3106.7646
$ws.Range("I132").Value = 2509
$ws.Range("K132").Value = 7527
$ws.Range("M132").Value = -4997

$ws.Range("H136").Value = 1886.2354
$ws.Range("I136").Value = 1632.25
$ws.Range("J136").Value = 5950
$ws.Range("K136").Value = 4896.75
$ws.Range("L136").Value = 17850
$ws.Range("M136").Value = -2346.75
$ws.Range("N136").Value = -22950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2245.7886
$ws.Range("J121").Value = 2360.2449
$ws.Range("L121").Value = 7080.734700000001
$ws.Range("N121").Value = -9700.734700000001

$ws.Range("H131").Value = 14712979
$ws.Range("I131").Value = 100042500
$ws.Range("J131").Value = 993.0345
$ws.Range("K131").Value = 300127500
$ws.Range("L131").Value = 2979.1035
$ws.Range("M131").Value = -300122460
$ws.Range("N131").Value = -13059.1035

$ws.Range("H132").Value = 5049.8335
$ws.Range("I132").Value = 1570.2858
$ws.Range("J132").Value = 7264.091
$ws.Range("K132").Value = 14132.5722
$ws.Range("L132").Value = 65376.819
$ws.Range("M132").Value = -11602.5722
$ws.Range("N132").Value = -70436.819

$ws.Range("H136").Value = 3003.6875
$ws.Range("I136").Value = 2466
$ws.Range("J136").Value = 3899.8333
$ws.Range("K136").Value = 7398
$ws.Range("L136").Value = 11699.4999
$ws.Range("M136").Value = -2298
$ws.Range("N136").Value = -21899.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11188.5
$ws.Range("J5").Value = 11297
$ws.Range("L5").Value = 11297
$ws.Range("N5").Value = -11521

$ws.Range("H46").Value = 39730.332
$ws.Range("J46").Value = 39730.332
$ws.Range("L46").Value = 39730.332
$ws.Range("N46").Value = -40042.332

$ws.Range("H107").Value = 4274000
$ws.Range("I107").Value = 202.07143
$ws.Range("K107").Value = 202.07143
$ws.Range("M107").Value = 1717.92857

$ws.Range("H137").Value = 48750
$ws.Range("J137").Value = 48750
$ws.Range("L137").Value = 48750
$ws.Range("N137").Value = -58950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 477.93332
$ws.Range("I16").Value = 485.13794
$ws.Range("J16").Value = 269
$ws.Range("K16").Value = 485.13794
$ws.Range("L16").Value = 269
$ws.Range("M16").Value = -315.13794
$ws.Range("N16").Value = -609

$ws.Range("H68").Value = 842.95123
$ws.Range("I68").Value = 707.20514
$ws.Range("J68").Value = 3490
$ws.Range("K68").Value = 707.20514
$ws.Range("L68").Value = 3490
$ws.Range("M68").Value = 41.79485999999997
$ws.Range("N68").Value = -4988

$ws.Range("H69").Value = 520000
$ws.Range("J69").Value = 520000
$ws.Range("L69").Value = 520000
$ws.Range("N69").Value = -521622

$ws.Range("H71").Value = 842.95123
$ws.Range("I71").Value = 707.20514
$ws.Range("J71").Value = 3490
$ws.Range("K71").Value = 3536.0257
$ws.Range("L71").Value = 17450
$ws.Range("M71").Value = 207.9742999999999
$ws.Range("N71").Value = -24938

$ws.Range("H72").Value = 520000
$ws.Range("J72").Value = 520000
$ws.Range("L72").Value = 1560000
$ws.Range("N72").Value = -1568112

$ws.Range("H110").Value = 39800
$ws.Range("J110").Value = 39800
$ws.Range("L110").Value = 39800
$ws.Range("N110").Value = -47980

$ws.Range("H122").Value = 6470.5293
$ws.Range("I122").Value = 4955.4443
$ws.Range("J122").Value = 8175
$ws.Range("K122").Value = 14866.3329
$ws.Range("L122").Value = 24525
$ws.Range("M122").Value = -12416.3329
$ws.Range("N122").Value = -29425

$ws.Range("H136").Value = 2773.1
$ws.Range("I136").Value = 1426.9546
$ws.Range("K136").Value = 4280.8638
$ws.Range("M136").Value = -1730.8638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9808948
$ws.Range("I132").Value = 6653.5884
$ws.Range("J132").Value = 19611242
$ws.Range("K132").Value = 19960.7652
$ws.Range("L132").Value = 58833726
$ws.Range("M132").Value = -17430.7652
$ws.Range("N132").Value = -58838786

$ws.Range("H136").Value = 4514.9546
$ws.Range("I136").Value = 1315.4166
$ws.Range("K136").Value = 3946.2498
$ws.Range("M136").Value = -1396.2498
